$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.26"
$ws.Range("E2").Value = "'4.87%"

$ws.Range("D3").Value = "'27.63"
$ws.Range("E3").Value = "'-2.75%"

$ws.Range("D4").Value = "'5.229"
$ws.Range("E4").Value = "'-0.42%"

$ws.Range("D5").Value = "'0.05923"
$ws.Range("E5").Value = "'3.83%"

$ws.Range("D6").Value = "'6.693"
$ws.Range("E6").Value = "'1.03%"

$ws.Range("D7").Value = "'0.8678"
$ws.Range("E7").Value = "'2.01%"

$ws.Range("D8").Value = "'1.036"
$ws.Range("E8").Value = "'15.45%"

$ws.Range("E9").Value = "'3.68%"

$ws.Range("E10").Value = "'1.65%"

$ws.Range("D11").Value = "'0.03258"
$ws.Range("E11").Value = "'3.19%"

$ws.Range("D12").Value = "'0.09216"

$ws.Range("D13").Value = "'0.001542"
$ws.Range("E13").Value = "'1.03%"

$ws.Range("D14").Value = "'0.0006084"
$ws.Range("E14").Value = "'2.11%"

$ws.Range("D15").Value = "'0.005798"
$ws.Range("E15").Value = "'-2.13%"

$ws.Range("D16").Value = "'3.483"
$ws.Range("E16").Value = "'-0.20%"

$ws.Range("D17").Value = "'3.272"
$ws.Range("E17").Value = "'2.22%"

$ws.Range("D18").Value = "'2.205"
$ws.Range("E18").Value = "'1.39%"

$ws.Range("D19").Value = "'0.3150"
$ws.Range("E19").Value = "'-0.62%"

$ws.Range("D20").Value = "'0.03593"
$ws.Range("E20").Value = "'9.40%"

$ws.Range("D21").Value = "'0.1308"
$ws.Range("E21").Value = "'2.48%"

$ws.Range("D22").Value = "'3.534"
$ws.Range("E22").Value = "'0.76%"

$ws.Range("D23").Value = "'0.04174"
$ws.Range("E23").Value = "'2.32%"

$ws.Range("E24").Value = "'1.55%"

$ws.Range("D25").Value = "'0.001220"
$ws.Range("E25").Value = "'-0.11%"

$ws.Range("D26").Value = "'0.004526"
$ws.Range("E26").Value = "'8.99%"

$ws.Range("E27").Value = "'0.12%"

$ws.Range("D28").Value = "'0.0001939"
$ws.Range("E28").Value = "'33.88%"

$ws.Range("E40").Value = "'1.11%"

$ws.Range("D41").Value = "'0.005453"
$ws.Range("E41").Value = "'4.45%"

$ws.Range("D42").Value = "'0.1106"
$ws.Range("E42").Value = "'3.80%"

$ws.Range("D43").Value = "'0.002461"
$ws.Range("E43").Value = "'11.95%"

$ws.Range("D44").Value = "'0.009840"
$ws.Range("E44").Value = "'7.54%"

$ws.Range("D45").Value = "'0.00005438"
$ws.Range("E45").Value = "'3.28%"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.11%"

$ws.Range("E47").Value = "'4.00%"

$ws.Range("E48").Value = "'-4.76%"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.11%"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.11%"
